$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '
              •
              COVID-19 in children of Xiangyang city is often family acquired and not serious, with favorable outcomes.
'
$ws.Range('E2').Value = '[Wenliang%Song%songwl2002jason@163.com%0, Junhua%Li%m13797631535@163.com%0, Ning%Zou%xdzning0215@163.com%1, Wenhe%Guan%guanwenhe@126.com%1, Jiali%Pan%cmu107819@163.com%1, Wei%Xu%tomxu.123@163.com%0]'
$ws.Range('I2').Value = ''
$ws.Range('J2').Value = 'Elsevier B.V.'
$ws.Range('D3').Value = 'We reported the clinical characteristics of a case series of 10 patients with coronavirus disease 2019 (COVID‐19) aged from 1 year to 18 years.
 Seven patients had contact with confirmed COVID‐19 family members before onset.
 Fever (4 [40.0%]) and cough (3 [30.0%]) were the most common symptoms.
 No patient showed leucopenia and lymphopenia on admission.
 Pneumonia was observed in chest CT images in 5 (50.0%) patients.
 Five (50.0%) patients received antiviral treatment.
 No patient had severe complications or developed a severe illness in our study.
 Our study indicated that COVID‐19 children present less severe symptoms and have better outcomes.
'
$ws.Range('E3').Value = '[Li%Zhu%NULL%0, Jian%Wang%NULL%0, Rui%Huang%NULL%0, Longgen%Liu%NULL%0, Haiyan%Zhao%NULL%0, Chao%Wu%dr.wu@nju.edu.cn%0, Chuanwu%Zhu%zhuchw@126.com%0, Chuanwu%Zhu%zhuchw@126.com%0]'
$ws.Range('I3').Value = ''
$ws.Range('J3').Value = 'John Wiley and Sons Inc.'
$ws.Range('D4').Value = 'We present clinical and chest computed tomography (CT) features of 5 cases of pediatric patients with 2019 novel coronavirus.
 Two patients had fever and dry cough, whereas the rest of 3 patients were asymptomatic.
 Three patients had unilateral ground glass opacities with or without consolidation in the subpleural region on high-resolution chest CT, 1 patient had bilateral ground glass opacities, and 1 patient was negative for CT.
 We note that up to 66.7% asymptomatic patients had pulmonary lesions, so the asymptomatic children with Wuhan contact are recommended to do a 2019 novel coronavirus real-time fluorescence polymerase chain reaction screening.
 Unlike adult patients, only a small amount of patients had multilobes affected, so we speculate that the pediatric patients generally have milder CT findings than adults.
'
$ws.Range('E4').Value = '[Mengqi%Liu%NULL%0, Zongbiao%Song%NULL%1, Kaihu%Xiao%NULL%0]'
$ws.Range('I4').Value = ''
$ws.Range('J4').Value = 'Lippincott Williams &amp; Wilkins'
$ws.Range('D5').Value = 'Background
id="Par1">Recently, the World Health Organization has declared the coronavirus disease 2019 (COVID-19) outbreak a public health emergency of international concern.

 So far, however, limited data are available for children.

 Therefore, we aimed to investigate the clinical and chest CT imaging characteristics of COVID-19 in preschool children.


Methods
id="Par2">From January 26, 2020 to February 20, 2020, the clinical and initial chest CT imaging data of eight preschool children with laboratory-confirmed COVID-19 from two hospitals were retrospectively collected.

 The chest CT imaging characteristics, including the distribution, shape, and density of lesions, and the pleural effusion, pleural changes, and enlarged lymph nodes were evaluated.


Results
id="Par3">Two cases (25%) were classified as mild type, and they showed no obvious abnormal CT findings or minimal pleural thickening on the right side.

 Five cases (62.5%) were classified as moderate type.

 Among these patients, one case showed consolidation located in the subpleural region of the right upper lobe, with thickening in the adjacent pleura; one case showed multiple consolidation and ground-glass opacities with blurry margins; one case displayed bronchial pneumonia-like changes in the left upper lobe; and two cases displayed asthmatic bronchitis-like changes.

 One case (12.5%) was classified as critical type and showed bronchial pneumonia-like changes in the bilateral lungs, presenting blurred and messy bilateral lung markings and multiple patchy shadows scattered along the lung markings with blurry margins.


Conclusions
id="Par4">The chest CT findings of COVID-19 in preschool children are atypical and various.

 Accurate diagnosis requires a comprehensive evaluation of epidemiological, clinical, laboratory and CT imaging data.


'
$ws.Range('E5').Value = '[Yang%Li%NULL%0, Jianghui%Cao%NULL%1, Xiaolong%Zhang%NULL%1, Guangzhi%Liu%NULL%1, Xiaxia%Wu%NULL%1, Baolin%Wu%WBLlin0129@163.com%1]'
$ws.Range('I5').Value = ''
$ws.Range('J5').Value = 'BioMed Central'
$ws.Range('C6').Value = 'Unknown Title'
$ws.Range('E6').Value = '[]'
$ws.Range('F6').Value = 'not found'
$ws.Range('G6').Value = 'N/A'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '1970-01-01'
$ws.Range('H6').ClearFormats()
$ws.Range('J6').Value = ''
$ws.Range('E7').Value = '[Y.%Lu%NULL%0, H.%Wen%NULL%1, D.%Rong%NULL%1, Z.%Zhou%NULL%1, H.%Liu%NULL%1]'
$ws.Range('I7').Value = ''
$ws.Range('J7').Value = 'The Royal College of Radiologists. Published by Elsevier Ltd.'
$ws.Range('E8').Value = '[Lan%Lan%NULL%0, Dan%Xu%NULL%0, Dan%Xu%NULL%0, Chen%Xia%NULL%2, Chen%Xia%NULL%0, Shaokang%Wang%NULL%2, Shaokang%Wang%NULL%0, Minhua%Yu%NULL%2, Minhua%Yu%NULL%0, Haibo%Xu%NULL%0, Haibo%Xu%NULL%0]'
$ws.Range('I8').Value = ''
$ws.Range('J8').Value = 'The Korean Society of Radiology'
$ws.Range('E9').Value = '[Jafar%Soltani%NULL%0, Iraj%Sedighi%NULL%1, Zohreh%Shalchi%NULL%1, Ghazal%Sami%NULL%1, Borhan%Moradveisi%NULL%1, Soheila%Nahidi%NULL%1]'
$ws.Range('I9').Value = ''
$ws.Range('J9').Value = 'Kare Publishing'

$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()
$ws.Rows(5).AutoFit()
